# Updated cryptos list data (Price / Volume(1h)) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.690.48"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.035.54"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.89"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.95"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "3.513.43"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.51"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.74"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "3.036.56"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.50"
$ws.Range("E18").Value = "  -15.06%  "
$ws.Range("D19").Value = "51.718.32"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.09"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.47"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.77"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("E25").Value = "  -7.59%  "
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("E27").Value = "  +11.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  +4.85%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.24"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.30"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.07"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.56"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.82"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +4.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.301"
$ws.Range("E39").Value = "  +18.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.03"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "127.31"
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("E45").Value = "  +6.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.69"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("D49").Value = "2.028.98"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "3.336.07"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  +1.33%  "
